# 自动更新价格数据: insert a new most-recent-date row (2025-11-30) at the
# top of the data (row 2), pushing the existing date rows down by one.
# The values for the new row repeat the same commodity prices as the
# rest of the table (783.5 / 1112 / 3610), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2..11 down to 3..12 by inserting a fresh row at row 2.
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting copied from the row above (bold
# header style); strip it so the new data row matches the plain style
# used by every other data row.
$ws.Range("A2:D2").ClearFormats()

# Force the date cell to be stored as text (matching the other date
# cells, which are plain text like "2025-11-29" rather than real Excel
# dates), then clear the number-format style that setting "@" added so
# the cell ends up with the same default style as its neighbours.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-11-30"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
$ws.Range("A2:D2").ClearFormats()
